# Applies the fidi.xlsx watchlist update:
#  - refreshes many of the ticker symbols in columns B, C, E, F for rows 2-26
#  - extends the table with 16 new rows (27-42) carrying columns A (index) and C (ticker)
#  - grows the used range from A1:F26 to A1:F42

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Updates to existing rows (2-26) -------------------------------------
$ws.Range("B2").Value = "NSE:HIRECT"
$ws.Range("C2").Value = "NSE:ACC"
$ws.Range("E2").Value = "NSE:APLAPOLLO"
$ws.Range("F2").Value = "NSE:LAURUSLABS"

$ws.Range("B3").Value = "NSE:LAURUSLABS"
$ws.Range("C3").Value = "NSE:AIRAN"
$ws.Range("E3").Value = "NSE:BHARATFORG"

$ws.Range("B4").Value = "NSE:MON100"
$ws.Range("C4").Value = "NSE:ASHAPURMIN"
$ws.Range("E4").Value = "NSE:CDSL"

$ws.Range("B5").Value = "NSE:PHARMABEES"
$ws.Range("C5").Value = "NSE:AURIONPRO"
$ws.Range("E5").Value = "NSE:DMART"

$ws.Range("B6").Value = "NSE:PUNJABCHEM"
$ws.Range("C6").Value = "NSE:BBTC"
$ws.Range("E6").Value = "NSE:GAIL"

$ws.Range("B7").Value = ""
$ws.Range("C7").Value = "NSE:BIKAJI"
$ws.Range("E7").Value = "NSE:INDUSINDBK"

$ws.Range("B8").Value = ""
$ws.Range("C8").Value = "NSE:BVCL"
$ws.Range("E8").Value = "NSE:IRCTC"

$ws.Range("B9").Value = ""
$ws.Range("C9").Value = "NSE:CARYSIL"
$ws.Range("E9").Value = "NSE:LT"

$ws.Range("C10").Value = "NSE:CASTROLIND"
$ws.Range("E10").Value = "NSE:ONGC"

$ws.Range("C11").Value = "NSE:CGPOWER"
$ws.Range("C12").Value = "NSE:CHEMPLASTS"
$ws.Range("C13").Value = "NSE:CONCORDBIO"
$ws.Range("C14").Value = "NSE:CUMMINSIND"
$ws.Range("C15").Value = "NSE:DCBBANK"
$ws.Range("C16").Value = "NSE:DCMSHRIRAM"
$ws.Range("C17").Value = "NSE:DMCC"
$ws.Range("C18").Value = "NSE:FINCABLES"
$ws.Range("C19").Value = "NSE:GODFRYPHLP"
$ws.Range("C20").Value = "NSE:GTPL"
$ws.Range("C21").Value = "NSE:HGS"
$ws.Range("C22").Value = "NSE:INTLCONV"
$ws.Range("C23").Value = "NSE:IOLCP"
$ws.Range("C24").Value = "NSE:IONEXCHANG"
$ws.Range("C25").Value = "NSE:LATENTVIEW"
$ws.Range("C26").Value = "NSE:LOVABLE"

# --- New rows 27-42 ---------------------------------------------------------
# Column A keeps the same bold/bordered "index" style as the rows above it,
# so copy the formatting from the last existing row before filling values in.
$ws.Range("A26").Copy()
$ws.Range("A27:A42").PasteSpecial(-4122)

$ws.Cells.Item(27, 1).Value = 25
$ws.Range("C27").Value = "NSE:MAITHANALL"

$ws.Cells.Item(28, 1).Value = 26
$ws.Range("C28").Value = "NSE:MANAKALUCO"

$ws.Cells.Item(29, 1).Value = 27
$ws.Range("C29").Value = "NSE:MGL"

$ws.Cells.Item(30, 1).Value = 28
$ws.Range("C30").Value = "NSE:MOIL"

$ws.Cells.Item(31, 1).Value = 29
$ws.Range("C31").Value = "NSE:MONARCH"

$ws.Cells.Item(32, 1).Value = 30
$ws.Range("C32").Value = "NSE:MOTILALOFS"

$ws.Cells.Item(33, 1).Value = 31
$ws.Range("C33").Value = "NSE:MRPL"

$ws.Cells.Item(34, 1).Value = 32
$ws.Range("C34").Value = "NSE:NATCOPHARM"

$ws.Cells.Item(35, 1).Value = 33
$ws.Range("C35").Value = "NSE:NESTLEIND"

$ws.Cells.Item(36, 1).Value = 34
$ws.Range("C36").Value = "NSE:NTPC"

$ws.Cells.Item(37, 1).Value = 35
$ws.Range("C37").Value = "NSE:OBCL"

$ws.Cells.Item(38, 1).Value = 36
$ws.Range("C38").Value = "NSE:ORIENTCER"

$ws.Cells.Item(39, 1).Value = 37
$ws.Range("C39").Value = "NSE:PERSISTENT"

$ws.Cells.Item(40, 1).Value = 38
$ws.Range("C40").Value = "NSE:PVP"

$ws.Cells.Item(41, 1).Value = 39
$ws.Range("C41").Value = "NSE:REFEX"

$ws.Cells.Item(42, 1).Value = 40
$ws.Range("C42").Value = "NSE:RELIANCE"

Write-Output "fidi.xlsx watchlist updated: rows 2-26 refreshed, rows 27-42 added"
